$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8568683862686157
$ws.Range("B1").Value = 2.993264198303223
$ws.Range("C1").Value = 3.076652765274048
$ws.Range("D1").Value = 1.744043707847595
$ws.Range("E1").Value = 1.339290857315063
